$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198, pushing the existing rows 198:210 down to 199:211.
$ws.Rows("198:198").Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(198, 1).Value = 11
$ws.Cells.Item(198, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(198, 3).Value = "Bíobío"
$ws.Cells.Item(198, 4).Value = 45106
$ws.Cells.Item(198, 4).Style = $ws.Cells.Item(199, 4).Style
$ws.Cells.Item(198, 4).NumberFormat = $ws.Cells.Item(199, 4).NumberFormat
$ws.Cells.Item(198, 5).Value = 8
$ws.Cells.Item(198, 6).Value = 100112021
$ws.Cells.Item(198, 7).Value = "Ají"
$ws.Cells.Item(198, 8).Value = "Americana (o)"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 40
$ws.Cells.Item(198, 11).Value = 28000
$ws.Cells.Item(198, 12).Value = 30000
$ws.Cells.Item(198, 13).Value = 29000
$ws.Cells.Item(198, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(198, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(198, 16).Value = 1160
$ws.Cells.Item(198, 17).Value = 25
$ws.Cells.Item(198, 18).Value = "Hortaliza"
